$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Cells.Item(2, 4).Value = "43.272.61"
$ws.Cells.Item(2, 5).Value = "  -1.70%  "
$ws.Cells.Item(3, 4).Value = "2.337.14"
$ws.Cells.Item(3, 5).Value = "  +3.23%  "
$ws.Cells.Item(4, 5).Value = "  -0.07%  "
Set-TextValue $ws.Cells.Item(5, 4) "0.649"
$ws.Cells.Item(5, 5).Value = "  +1.68%  "
Set-TextValue $ws.Cells.Item(6, 4) "230.79"
Set-TextValue $ws.Cells.Item(7, 4) "65.04"
$ws.Cells.Item(7, 5).Value = "  +0.65%  "
$ws.Cells.Item(8, 5).Value = "  -0.05%  "
Set-TextValue $ws.Cells.Item(9, 4) "0.453"
$ws.Cells.Item(9, 5).Value = "  +0.85%  "
Set-TextValue $ws.Cells.Item(10, 4) "0.0949"
$ws.Cells.Item(10, 5).Value = "  -5.05%  "
Set-TextValue $ws.Cells.Item(11, 4) "56.78"
$ws.Cells.Item(11, 5).Value = "  -0.80%  "
Set-TextValue $ws.Cells.Item(12, 4) "26.56"
$ws.Cells.Item(12, 5).Value = "  -3.60%  "
$ws.Cells.Item(13, 4).Value = "2.679.65"
$ws.Cells.Item(13, 5).Value = "  +2.90%  "
Set-TextValue $ws.Cells.Item(15, 4) "15.27"
$ws.Cells.Item(15, 5).Value = "  -3.20%  "
Set-TextValue $ws.Cells.Item(16, 4) "6.22"
$ws.Cells.Item(16, 5).Value = "  +1.80%  "
Set-TextValue $ws.Cells.Item(17, 4) "0.837"
$ws.Cells.Item(17, 5).Value = "  -0.37%  "
$ws.Cells.Item(18, 4).Value = "2.333.48"
$ws.Cells.Item(18, 5).Value = "  +2.15%  "
$ws.Cells.Item(19, 4).Value = "43.149.76"
$ws.Cells.Item(19, 5).Value = "  -1.95%  "
$ws.Cells.Item(20, 4).Value = "0.0₃0971"
$ws.Cells.Item(20, 5).Value = "  -3.61%  "
Set-TextValue $ws.Cells.Item(21, 4) "73.56"
$ws.Cells.Item(21, 5).Value = "  -0.39%  "
$ws.Cells.Item(22, 5).Value = "  +0.42%  "
Set-TextValue $ws.Cells.Item(23, 4) "247.60"
$ws.Cells.Item(23, 5).Value = "  -2.28%  "
Set-TextValue $ws.Cells.Item(24, 4) "3.90"
$ws.Cells.Item(24, 5).Value = "  +18.41%  "
$ws.Cells.Item(25, 5).Value = "  -0.05%  "
$ws.Cells.Item(26, 5).Value = "  -1.41%  "
$ws.Cells.Item(27, 5).Value = "  -1.65%  "
Set-TextValue $ws.Cells.Item(28, 4) "9.82"
$ws.Cells.Item(28, 5).Value = "  -3.04%  "
Set-TextValue $ws.Cells.Item(29, 4) "174.69"
$ws.Cells.Item(29, 5).Value = "  +1.77%  "
Set-TextValue $ws.Cells.Item(30, 4) "22.16"
$ws.Cells.Item(30, 5).Value = "  +5.77%  "
$ws.Cells.Item(31, 5).Value = "  +4.61%  "
$ws.Cells.Item(32, 5).Value = "  -8.77%  "
$ws.Cells.Item(33, 5).Value = "  +0.06%  "
Set-TextValue $ws.Cells.Item(34, 4) "5.00"
$ws.Cells.Item(34, 5).Value = "  +3.71%  "
Set-TextValue $ws.Cells.Item(35, 4) "0.0682"
$ws.Cells.Item(35, 5).Value = "  -3.48%  "
Set-TextValue $ws.Cells.Item(36, 4) "4.94"
$ws.Cells.Item(36, 5).Value = "  +0.95%  "
$ws.Cells.Item(37, 5).Value = "  +5.95%  "
Set-TextValue $ws.Cells.Item(38, 4) "6.44"
$ws.Cells.Item(38, 5).Value = "  -1.13%  "
$ws.Cells.Item(39, 5).Value = "  -6.62%  "
Set-TextValue $ws.Cells.Item(40, 4) "0.0250"
$ws.Cells.Item(40, 5).Value = "  -4.21%  "
$ws.Cells.Item(41, 5).Value = "  -0.06%  "
Set-TextValue $ws.Cells.Item(42, 4) "8.86"
$ws.Cells.Item(42, 5).Value = "  +7.46%  "
$ws.Cells.Item(43, 5).Value = "  +1.27%  "
$ws.Cells.Item(44, 5).Value = "  +6.05%  "
Set-TextValue $ws.Cells.Item(45, 4) "98.17"
$ws.Cells.Item(45, 5).Value = "  -0.32%  "
$ws.Cells.Item(46, 5).Value = "  -1.09%  "
Set-TextValue $ws.Cells.Item(47, 4) "4.36"
$ws.Cells.Item(47, 5).Value = "  -0.49%  "
Set-TextValue $ws.Cells.Item(48, 4) "0.0941"
$ws.Cells.Item(48, 5).Value = "  -4.66%  "
$ws.Cells.Item(49, 4).Value = "1.433.41"
$ws.Cells.Item(49, 5).Value = "  -1.06%  "
Set-TextValue $ws.Cells.Item(50, 4) "9.81"
$ws.Cells.Item(50, 5).Value = "  -7.96%  "
$ws.Cells.Item(51, 2).Value = "TerraClassic"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
Set-TextValue $ws.Cells.Item(51, 4) "0.000202"
$ws.Cells.Item(51, 5).Value = "  -9.87%  "
